$wb = $excel.ActiveWorkbook
$sys = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) Insert a brand-new column before column Y. This pushes the existing
#    "web" / "webalert" / "webcookie" / "ws" / "ws.async" / "xml" lists one
#    column to the right (Y->Z, Z->AA, AA->AB, AB->AC, AC->AD, AD->AE) for
#    every row (1-129), and leaves everything left of Y untouched.
# ---------------------------------------------------------------------------
$sys.Columns("Y:Y").Insert()

# Populate the freshly inserted column with the new "text" category: a
# header in row 1 and its single command in row 2.
$sys.Cells.Item(1, 25).Value = "text"
$sys.Cells.Item(2, 25).Value = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------------
# 2) "target" list (column A): insert a new "text" entry alphabetically
#    between "step" and "web" (row 25), pushing A25:A30 down to A26:A31.
#    Only column A is affected - every other column on those rows keeps its
#    original value, so we shift cell-by-cell instead of using a row/range
#    Insert (which would move entire rows).
# ---------------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $sys.Cells.Item($r + 1, 1).Value = $sys.Cells.Item($r, 1).Text
}
$sys.Cells.Item(25, 1).Value = "text"

# ---------------------------------------------------------------------------
# 3) "base" list (column E): insert the new "outputToCloud(resource)"
#    command alphabetically between "macro(file,sheet,name)" and
#    "prependText(var,prependWith)" (row 22), pushing E22:E38 down to
#    E23:E39. Again, only column E moves.
# ---------------------------------------------------------------------------
for ($r = 38; $r -ge 22; $r--) {
    $sys.Cells.Item($r + 1, 5).Value = $sys.Cells.Item($r, 5).Text
}
$sys.Cells.Item(22, 5).Value = "outputToCloud(resource)"

# ---------------------------------------------------------------------------
# 4) Update the defined names so they point at the resized/shifted ranges,
#    and register the new "text" name.
# ---------------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
